$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.06754433333333333"
$ws.Range("H2").Value = [double]"0.202633"
$ws.Range("I2").Value = [double]"0.006855017925354449"
$ws.Range("J2").Value = [double]"0.006855017925354449"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.3333333333333333"
$ws.Range("M2").Value = [double]"0.2189473333333334"
$ws.Range("N2").Value = [double]"0.656842"
$ws.Range("O2").Value = [double]"0.009402596261870986"
$ws.Range("P2").Value = [double]"0.009402596261870984"
$ws.Range("Q2").Value = [double]"0.01478865166511111"
$ws.Range("R2").Value = [double]"0.133097864986"
$ws.Range("S2").Value = [double]"6.445496591999635E-05"
$ws.Range("T2").Value = [double]"6.445496591999633E-05"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.06754433333333333"
$ws.Range("H3").Value = [double]"0.202633"
$ws.Range("I3").Value = [double]"0.006855017925354449"
$ws.Range("J3").Value = [double]"0.006855017925354449"
$ws.Range("O3").Value = [double]"0.8622887582286424"
$ws.Range("P3").Value = [double]"0.8622887582286423"
$ws.Range("Q3").Value = [double]"1.356230526657445"
$ws.Range("R3").Value = [double]"12.206074739917"
$ws.Range("S3").Value = [double]"0.005911004894488972"
$ws.Range("T3").Value = [double]"0.005911004894488972"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.06754433333333333"
$ws.Range("H4").Value = [double]"0.202633"
$ws.Range("I4").Value = [double]"0.006855017925354449"
$ws.Range("J4").Value = [double]"0.006855017925354449"
$ws.Range("M4").Value = [double]"2.823530666666667"
$ws.Range("N4").Value = [double]"8.470592"
$ws.Range("O4").Value = [double]"0.1212552739852724"
$ws.Range("P4").Value = [double]"0.1212552739852723"
$ws.Range("Q4").Value = [double]"0.1907134965262222"
$ws.Range("R4").Value = [double]"1.716421468736"
$ws.Range("S4").Value = [double]"0.000831207076712807"
$ws.Range("T4").Value = [double]"0.0008312070767128069"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"0.3333333333333333"
$ws.Range("G5").Value = [double]"0.06754433333333333"
$ws.Range("H5").Value = [double]"0.202633"
$ws.Range("I5").Value = [double]"0.006855017925354449"
$ws.Range("J5").Value = [double]"0.006855017925354449"
$ws.Range("M5").Value = [double]"0.1642436666666667"
$ws.Range("N5").Value = [double]"0.492731"
$ws.Range("O5").Value = [double]"0.007053371524214274"
$ws.Range("P5").Value = [double]"0.007053371524214274"
$ws.Range("Q5").Value = [double]"0.01109372896922222"
$ws.Range("R5").Value = [double]"0.099843560723"
$ws.Range("S5").Value = [double]"4.835098823267348E-05"
$ws.Range("T5").Value = [double]"4.835098823267348E-05"
$ws.Range("I6").Value = [double]"0.7774992501642265"
$ws.Range("J6").Value = [double]"0.7774992501642265"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.3333333333333333"
$ws.Range("M6").Value = [double]"0.2189473333333334"
$ws.Range("N6").Value = [double]"0.656842"
$ws.Range("O6").Value = [double]"0.009402596261870986"
$ws.Range("P6").Value = [double]"0.009402596261870984"
$ws.Range("Q6").Value = [double]"1.677335596459334"
$ws.Range("R6").Value = [double]"15.096020368134"
$ws.Range("S6").Value = [double]"0.007310511543201651"
$ws.Range("T6").Value = [double]"0.007310511543201649"
$ws.Range("I7").Value = [double]"0.7774992501642265"
$ws.Range("J7").Value = [double]"0.7774992501642265"
$ws.Range("O7").Value = [double]"0.8622887582286424"
$ws.Range("P7").Value = [double]"0.8622887582286423"
$ws.Range("S7").Value = [double]"0.6704288629478115"
$ws.Range("T7").Value = [double]"0.6704288629478113"
$ws.Range("I8").Value = [double]"0.7774992501642265"
$ws.Range("J8").Value = [double]"0.7774992501642265"
$ws.Range("M8").Value = [double]"2.823530666666667"
$ws.Range("N8").Value = [double]"8.470592"
$ws.Range("O8").Value = [double]"0.1212552739852724"
$ws.Range("P8").Value = [double]"0.1212552739852723"
$ws.Range("Q8").Value = [double]"21.63081149604267"
$ws.Range("R8").Value = [double]"194.677303464384"
$ws.Range("S8").Value = [double]"0.0942758846020071"
$ws.Range("T8").Value = [double]"0.09427588460200709"
$ws.Range("I9").Value = [double]"0.7774992501642265"
$ws.Range("J9").Value = [double]"0.7774992501642265"
$ws.Range("M9").Value = [double]"0.1642436666666667"
$ws.Range("N9").Value = [double]"0.492731"
$ws.Range("O9").Value = [double]"0.007053371524214274"
$ws.Range("P9").Value = [double]"0.007053371524214274"
$ws.Range("Q9").Value = [double]"1.258255784159667"
$ws.Range("R9").Value = [double]"11.324302057437"
$ws.Range("S9").Value = [double]"0.005483991071206305"
$ws.Range("T9").Value = [double]"0.005483991071206305"
$ws.Range("G10").Value = [double]"1.941983333333333"
$ws.Range("H10").Value = [double]"5.825949999999999"
$ws.Range("I10").Value = [double]"0.1970902650714284"
$ws.Range("J10").Value = [double]"0.1970902650714283"
$ws.Range("K10").Value = [double]"1"
$ws.Range("L10").Value = [double]"0.3333333333333333"
$ws.Range("M10").Value = [double]"0.2189473333333334"
$ws.Range("N10").Value = [double]"0.656842"
$ws.Range("O10").Value = [double]"0.009402596261870986"
$ws.Range("P10").Value = [double]"0.009402596261870984"
$ws.Range("Q10").Value = [double]"0.425192072211111"
$ws.Range("R10").Value = [double]"3.826728649899999"
$ws.Range("S10").Value = [double]"0.001853160189611774"
$ws.Range("T10").Value = [double]"0.001853160189611774"
$ws.Range("G11").Value = [double]"1.941983333333333"
$ws.Range("H11").Value = [double]"5.825949999999999"
$ws.Range("I11").Value = [double]"0.1970902650714284"
$ws.Range("J11").Value = [double]"0.1970902650714283"
$ws.Range("O11").Value = [double]"0.8622887582286424"
$ws.Range("P11").Value = [double]"0.8622887582286423"
$ws.Range("Q11").Value = [double]"38.99330926739444"
$ws.Range("R11").Value = [double]"350.93978340655"
$ws.Range("S11").Value = [double]"0.1699487199273959"
$ws.Range("T11").Value = [double]"0.1699487199273959"
$ws.Range("G12").Value = [double]"1.941983333333333"
$ws.Range("H12").Value = [double]"5.825949999999999"
$ws.Range("I12").Value = [double]"0.1970902650714284"
$ws.Range("J12").Value = [double]"0.1970902650714283"
$ws.Range("M12").Value = [double]"2.823530666666667"
$ws.Range("N12").Value = [double]"8.470592"
$ws.Range("O12").Value = [double]"0.1212552739852724"
$ws.Range("P12").Value = [double]"0.1212552739852723"
$ws.Range("Q12").Value = [double]"5.483249495822221"
$ws.Range("R12").Value = [double]"49.34924546239999"
$ws.Range("S12").Value = [double]"0.023898234091066"
$ws.Range("T12").Value = [double]"0.023898234091066"
$ws.Range("G13").Value = [double]"1.941983333333333"
$ws.Range("H13").Value = [double]"5.825949999999999"
$ws.Range("I13").Value = [double]"0.1970902650714284"
$ws.Range("J13").Value = [double]"0.1970902650714283"
$ws.Range("M13").Value = [double]"0.1642436666666667"
$ws.Range("N13").Value = [double]"0.492731"
$ws.Range("O13").Value = [double]"0.007053371524214274"
$ws.Range("P13").Value = [double]"0.007053371524214274"
$ws.Range("Q13").Value = [double]"0.3189584632722222"
$ws.Range("R13").Value = [double]"2.870626169449999"
$ws.Range("S13").Value = [double]"0.001390150863354656"
$ws.Range("T13").Value = [double]"0.001390150863354656"
$ws.Range("E14").Value = [double]"2"
$ws.Range("F14").Value = [double]"0.6666666666666666"
$ws.Range("G14").Value = [double]"0.182832"
$ws.Range("H14").Value = [double]"0.548496"
$ws.Range("I14").Value = [double]"0.01855546683899075"
$ws.Range("J14").Value = [double]"0.01855546683899075"
$ws.Range("K14").Value = [double]"1"
$ws.Range("L14").Value = [double]"0.3333333333333333"
$ws.Range("M14").Value = [double]"0.2189473333333334"
$ws.Range("N14").Value = [double]"0.656842"
$ws.Range("O14").Value = [double]"0.009402596261870986"
$ws.Range("P14").Value = [double]"0.009402596261870984"
$ws.Range("Q14").Value = [double]"0.040030578848"
$ws.Range("R14").Value = [double]"0.360275209632"
$ws.Range("S14").Value = [double]"0.0001744695631375655"
$ws.Range("T14").Value = [double]"0.0001744695631375655"
$ws.Range("E15").Value = [double]"2"
$ws.Range("F15").Value = [double]"0.6666666666666666"
$ws.Range("G15").Value = [double]"0.182832"
$ws.Range("H15").Value = [double]"0.548496"
$ws.Range("I15").Value = [double]"0.01855546683899075"
$ws.Range("J15").Value = [double]"0.01855546683899075"
$ws.Range("O15").Value = [double]"0.8622887582286424"
$ws.Range("P15").Value = [double]"0.8622887582286423"
$ws.Range("Q15").Value = [double]"3.671104997456"
$ws.Range("R15").Value = [double]"33.039944977104"
$ws.Range("S15").Value = [double]"0.01600017045894609"
$ws.Range("T15").Value = [double]"0.01600017045894609"
$ws.Range("E16").Value = [double]"2"
$ws.Range("F16").Value = [double]"0.6666666666666666"
$ws.Range("G16").Value = [double]"0.182832"
$ws.Range("H16").Value = [double]"0.548496"
$ws.Range("I16").Value = [double]"0.01855546683899075"
$ws.Range("J16").Value = [double]"0.01855546683899075"
$ws.Range("M16").Value = [double]"2.823530666666667"
$ws.Range("N16").Value = [double]"8.470592"
$ws.Range("O16").Value = [double]"0.1212552739852724"
$ws.Range("P16").Value = [double]"0.1212552739852723"
$ws.Range("Q16").Value = [double]"0.516231758848"
$ws.Range("R16").Value = [double]"4.646085829632"
$ws.Range("S16").Value = [double]"0.00224994821548646"
$ws.Range("T16").Value = [double]"0.002249948215486459"
$ws.Range("E17").Value = [double]"2"
$ws.Range("F17").Value = [double]"0.6666666666666666"
$ws.Range("G17").Value = [double]"0.182832"
$ws.Range("H17").Value = [double]"0.548496"
$ws.Range("I17").Value = [double]"0.01855546683899075"
$ws.Range("J17").Value = [double]"0.01855546683899075"
$ws.Range("M17").Value = [double]"0.1642436666666667"
$ws.Range("N17").Value = [double]"0.492731"
$ws.Range("O17").Value = [double]"0.007053371524214274"
$ws.Range("P17").Value = [double]"0.007053371524214274"
$ws.Range("Q17").Value = [double]"0.030028998064"
$ws.Range("R17").Value = [double]"0.270260982576"
$ws.Range("S17").Value = [double]"0.0001308786014206396"
$ws.Range("T17").Value = [double]"0.0001308786014206396"
